# "Change set details added" - populate the "Change Sets" sheet with its
# first two data rows and make it the active/selected sheet & tab
# (previously "Frontend" was the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Sets")
$ws.Activate()

$ws.Range("A2").Value = "Page 1"
$ws.Range("A3").Value = "All pages"

# Leave the selection on the next empty row, as in the authored workbook.
[void]$ws.Range("A4").Select()
